$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Rename the column headers (Portuguese localization update) and make them bold,
# matching the already-bold "Data" header style.
$ws.Range("B1").Value = "Tipo de comentários"
$ws.Range("C1").Value = "Classificação de comentários"
$ws.Range("D1").Value = "Tipo de engajamento"
$ws.Range("E1").Value = "Classificação de engajamento"

$ws.Range("B1:E1").Font.Bold = $true
